$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as text so the date-like string isn't
# auto-converted into a date serial number (matches existing rows, which
# store the date as plain text).
$ws.Range("A34").NumberFormat = "@"

# Append the new mod-count data row for 2025/12/13
$ws.Range("A34").Value = "2025/12/13"
$ws.Range("B34").Value = "逃离鸭科夫"
$ws.Range("C34").Value = 1360

# Copy the formatting (centered alignment) from the previous data row so
# the new row matches the sheet's existing style.
$ws.Range("A33:C33").Copy()
$ws.Range("A34:C34").PasteSpecial(-4122)
